$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0002534235536586493
$ws.Range("A3").Value = 0.00011716933659045026
$ws.Range("H3").Value = 5.329891204833984
$ws.Range("A4").Value = 0.0001020246054395102
$ws.Range("H4").Value = 4.539682865142822
